$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.573.20"
$ws.Range("E2").Value = "  -0.78%  "
$ws.Range("D3").Value = "2.796.83"
$ws.Range("E3").Value = "  +0.21%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "355.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.06"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.554"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.10%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.625"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.83%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.80"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.22%  "
$ws.Range("E11").Value = "  +1.02%  "
$ws.Range("E12").Value = "  -1.68%  "
$ws.Range("E13").Value = "  +2.71%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.75"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.20%  "
$ws.Range("D15").Value = "3.235.68"
$ws.Range("E15").Value = "  +0.02%  "
$ws.Range("D16").Value = "2.796.78"
$ws.Range("E16").Value = "  +0.21%  "
$ws.Range("E17").Value = "  -1.03%  "
$ws.Range("D18").Value = "51.556.49"
$ws.Range("E18").Value = "  -0.75%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.73"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.78%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.15"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.48%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.32"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.28%  "
$ws.Range("E22").Value = "  -0.77%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.46"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.26%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "267.03"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.17%  "
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.96"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.24%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.164"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.32"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.17%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "37.14"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +8.25%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.24"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.64%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.23"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +8.60%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "52.28"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.23%  "
$ws.Range("E34").Value = "  +8.62%  "
$ws.Range("E35").Value = "  -6.11%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0853"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.94%  "
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.64"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.44%  "
$ws.Range("E39").Value = "  -2.72%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.98"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.74%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.115"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.07%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.48"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.86%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "119.85"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.48%  "
$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.19"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.65%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.70"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.58%  "
$ws.Range("D46").Value = "2.132.54"
$ws.Range("E46").Value = "  +1.99%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.40"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.15%  "
$ws.Range("E48").Value = "  +6.43%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.223"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +17.44%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.912"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.39%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.36"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +9.87%  "
